$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = 'maa://24702 (94.59), maa://25390 (96.21), maa://36681 (87.34)'
$ws.Range('L2').Value = '*maa://24633 (56.17), *maa://30515 (70.48), maa://39402 (92.42), *maa://34787 (72.73), ***maa://20792 (11.93), ***maa://29083 (27.78)'
$ws.Range('H3').Value = 'maa://21247 (98.61), *maa://22748 (60.0)'
$ws.Range('L3').Value = '*maa://22880 (65.0), maa://20276 (86.59), *maa://22749 (76.92)'
$ws.Range('O3').Value = '3'
$ws.Range('P3').Value = 'maa://21249 (94.58), maa://26254 (96.67), **maa://22738 (50.0)'
$ws.Range('T3').Value = 'maa://24617 (89.83), **maa://20790 (43.48), ***maa://37170 (16.42), maa://45854 (82.14)'
$ws.Range('X4').Value = '**maa://32495 (48.7), ***maa://31785 (22.22), maa://43217 (91.78), ***maa://36683 (28.26)'
$ws.Range('AF4').Value = '*maa://30062 (64.0), ***maa://26209 (13.04), *maa://39394 (65.38)'
$ws.Range('D5').Value = 'maa://21245 (84.36), maa://22744 (84.62)'
$ws.Range('X7').Value = 'maa://22399 (95.57), *maa://22758 (75.0)'
$ws.Range('A8').Value = '更新日期：2025.03.19 13:19:27'
$ws.Range('L9').Value = 'maa://22762 (92.39), maa://39552 (81.25)'
$ws.Range('AB9').Value = 'maa://28711 (87.1), ***maa://22740 (5.66), **maa://39938 (46.67), **maa://27377 (42.86), ***maa://25174 (19.05), *maa://45044 (66.67), maa://40166 (96.3)'
$ws.Range('AF9').Value = 'maa://26206 (88.1), *maa://22865 (51.85)'
$ws.Range('D10').Value = '***maa://25695 (18.52), ***maa://39951 (13.79), ***maa://34206 (19.23), ***maa://39243 (25.0), *maa://45271 (57.5)'
$ws.Range('X10').Value = 'maa://22301 (97.77), maa://45828 (86.67), maa://22726 (100.0)'
$ws.Range('AB11').Value = 'maa://29912 (97.33), maa://22516 (88.37), *maa://20794 (52.24)'
$ws.Range('H12').Value = 'maa://21867 (90.0), ***maa://45826 (25.0)'
$ws.Range('AB12').Value = 'maa://23669 (95.5), maa://36677 (93.94), maa://39872 (92.0)'
$ws.Range('D13').Value = 'maa://24999 (92.09), maa://36673 (92.21), maa://25001 (85.71)'
$ws.Range('AF13').Value = '**maa://22737 (34.25), maa://39883 (91.78), *maa://39885 (53.33)'
$ws.Range('H15').Value = 'maa://24304 (87.84), maa://21478 (89.19)'
$ws.Range('T16').Value = 'maa://22729 (94.55), *maa://28648 (69.57), maa://36674 (81.13)'
$ws.Range('AB16').Value = 'maa://26228 (95.1)'
$ws.Range('P17').Value = 'maa://23890 (80.95), *maa://24940 (67.86)'
$ws.Range('X18').Value = 'maa://21917 (96.94), maa://22741 (87.5)'
$ws.Range('L20').Value = 'maa://41331 (84.85)'
$ws.Range('X21').Value = 'maa://20110 (86.76), maa://34946 (91.11)'
$ws.Range('AF21').Value = 'maa://22524 (94.22), *maa://22432 (76.92)'
$ws.Range('D23').Value = '***maa://28036 (28.77), *maa://41753 (52.38)'
$ws.Range('L23').Value = 'maa://39756 (95.9), maa://39875 (94.52)'
$ws.Range('T23').Value = 'maa://24387 (82.05), maa://31212 (93.75)'
$ws.Range('D24').Value = '*maa://24368 (78.44), *maa://46650 (57.14)'
$ws.Range('X24').Value = 'maa://29988 (84.11), maa://23504 (93.33), **maa://22892 (40.14), *maa://25141 (77.1), *maa://36663 (77.5), ***maa://22815 (23.08)'
$ws.Range('D25').Value = 'maa://29753 (95.22)'
$ws.Range('AB26').Value = 'maa://42235 (94.59)'
$ws.Range('L27').Value = 'maa://28071 (90.48)'
$ws.Range('T28').Value = '*maa://29765 (64.37), maa://23263 (95.28)'
$ws.Range('X28').Value = 'maa://39929 (90.82), maa://41749 (91.4), ***maa://39723 (13.89)'
$ws.Range('AF29').Value = '*maa://24080 (68.85), maa://42865 (81.82), ***maa://34960 (8.33)'
$ws.Range('L30').Value = 'maa://30442 (95.38)'
$ws.Range('AB30').Value = 'maa://42979 (97.13), maa://45822 (100.0), *maa://45045 (80.0)'
$ws.Range('T32').Value = 'maa://42859 (96.12), maa://41108 (88.0), maa://41238 (97.2), maa://45523 (100.0)'
$ws.Range('P33').Value = 'maa://21956 (80.79), *maa://22730 (76.67)'
$ws.Range('O34').Value = '1'
$ws.Range('P34').Value = 'maa://48817 (100.0)'
$ws.Range('L37').Value = 'maa://45718 (98.17), *maa://47069 (73.33), maa://45789 (100.0)'
$ws.Range('H39').Value = 'maa://36670 (89.22), maa://25199 (84.82), maa://30434 (91.86), maa://45059 (85.0), ***maa://25036 (16.0), *maa://44165 (66.67)'
$ws.Range('T39').Value = 'maa://45788 (80.77), maa://47079 (92.0), *maa://45790 (73.33)'
$ws.Range('P41').Value = '**maa://35616 (40.0), maa://43177 (91.67)'
$ws.Range('H43').Value = 'maa://22525 (92.47), maa://21284 (85.71)'
$ws.Range('T45').Value = '**maa://39364 (38.89)'
$ws.Range('P49').Value = '*maa://39643 (63.64)'
$ws.Range('H53').Value = 'maa://32534 (94.2), **maa://32434 (33.33)'
$ws.Range('H62').Value = 'maa://42981 (95.35), maa://43903 (100.0)'
